$wb = $excel.ActiveWorkbook

# --- Sheet1: Week 14 (column Q) scores for rows 12,13,16,18,20,21,22,24,26,27,29,31 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("Q12").Value = 32
$ws1.Range("Q13").Value = 32
$ws1.Range("Q16").Value = 25
$ws1.Range("Q18").Value = 28
$ws1.Range("Q20").Value = 28
$ws1.Range("Q21").Value = 28
$ws1.Range("Q22").Value = 26
$ws1.Range("Q24").Value = 38
$ws1.Range("Q26").Value = 34
$ws1.Range("Q27").Value = 29
$ws1.Range("Q29").Value = 34
$ws1.Range("Q31").Value = 28

# --- THURSDAY SINGLES: Week 13 (column N) scores for rows 6,9,13,15 ---
$ws2 = $wb.Worksheets.Item("THURSDAY SINGLES")

$ws2.Range("N6").Value = 33
$ws2.Range("N9").Value = 28
$ws2.Range("N13").Value = 29
$ws2.Range("N15").Value = 27

# --- HANDICAPS: MICK SKINNER handicap bump ---
$ws3 = $wb.Worksheets.Item("HANDICAPS")

$ws3.Range("B3").Value = 6
